# Hoàn thiện Ngoại Trú
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# A2: patient/receipt id changed
$ws.Range("A2").Value = 3016

# New cell C2 inserted (Dob), shifting nothing else in the row structurally
$ws.Range("C2").Value = "2000-01-01T00:00:00+07:00"

# E2: insurance card number changed
$ws.Range("E2").Value = 46200608016

# M2 ("Address" column) cleared - no longer populated
$ws.Range("M2").ClearContents()

# AM2 / AN2 updated
$ws.Range("AM2").Value = 1
$ws.Range("AN2").Value = 0

# AQ2 ("InsCheckedMessage") cleared - no longer populated
$ws.Range("AQ2").ClearContents()

# Update selection to match new active cell
$ws.Range("AO2").Select() | Out-Null
